$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A63").Value = "Aste Andrea "
$ws.Range("B63").Value = "Matteo Pilati | Pinguini Trentini"
$ws.Range("C63").Value = "Andrea Roveda | Pinguini Trentini"
$ws.Range("D63").Value = "Alessio Bragagna | SHARK ATTACK"
$ws.Range("E63").Value = "Edoardo Pomarolli | Modium"
$ws.Range("F63").Value = "Matteo Maraner | GREP"
